$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new column O data for year 2021
$ws.Cells.Item(4, 15).Value = 2021
$ws.Cells.Item(5, 15).Value = 1.5020015556876996

# Copy styles from column N (14) to column O (15) for rows 4 and 5
$ws.Cells.Item(4, 14).Copy() | Out-Null
$ws.Cells.Item(4, 15).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(5, 14).Copy() | Out-Null
$ws.Cells.Item(5, 15).PasteSpecial(-4122) | Out-Null

# Update selection
$ws.Range("Q5").Select() | Out-Null
